# Add two new columns (K, L) with qkeras fixed-point accuracy results.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells - copy style from existing header cell J1 so the new
# headers look consistent (bold, bordered, centered) with the rest.
$ws.Range("K1").Value = "fxppo2_accuracy_qkeras"
$ws.Range("L1").Value = "orig-fxppo2-drop_qkeras"
$ws.Range("J1").Copy()
$ws.Range("K1:L1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("K1").Value = "fxppo2_accuracy_qkeras"
$ws.Range("L1").Value = "orig-fxppo2-drop_qkeras"

# Data rows 2-21 for column K (fxppo2_accuracy_qkeras)
$kValues = @(
    0.6071428571428571,
    0.6517857142857143,
    0.4553571428571428,
    0.5,
    0.4910714285714285,
    0.5,
    0.3035714285714285,
    0.6071428571428571,
    0.4553571428571428,
    0.4285714285714285,
    0.4285714285714285,
    0.5,
    0.5357142857142857,
    0.4821428571428572,
    0.4196428571428572,
    0.4910714285714285,
    0.5535714285714286,
    0.5446428571428571,
    0.3839285714285715,
    0.3035714285714285
)

# Data rows 2-21 for column L (orig-fxppo2-drop_qkeras)
$lValues = @(
    0.1785714285714286,
    0.1428571428571428,
    0.3125000000000001,
    0.0982142857142857,
    0.1071428571428572,
    0.2142857142857143,
    [double]"-5.551115123125783e-17",
    0.2053571428571429,
    -0.008928571428571341,
    0.08035714285714285,
    0.01785714285714296,
    0.1339285714285714,
    0.0267857142857143,
    0.01785714285714285,
    0.3303571428571428,
    0.008928571428571452,
    0.0446428571428571,
    0.2232142857142858,
    0.07142857142857134,
    [double]"-5.551115123125783e-17"
)

for ($i = 0; $i -lt 20; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 11).Value = $kValues[$i]
    $ws.Cells.Item($row, 12).Value = $lValues[$i]
}
